# Normalize the "Recorded By" column (G) so that the "System" entry
# (exact case) is moved to the front of the comma-separated list of
# recorders, preserving the relative order of the remaining entries.
# Special case: if the list already starts with a case-insensitive
# variant of "system" (but not exactly "System"), the "System" entry is
# inserted right after that leading entry instead of before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notmatch "System") { continue }

    $parts = @($val -split ", ")
    $hasSystem = $false
    $rest = @()
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        } else {
            $rest = $rest + ,$p
        }
    }

    if (-not $hasSystem) { continue }

    if ($rest.Count -gt 0 -and $rest[0].ToLower().Equals("system") -and -not $rest[0].Equals("System")) {
        $newParts = @($rest[0]) + @("System") + @($rest[1..($rest.Count - 1)])
    } else {
        $newParts = @("System") + $rest
    }

    $newVal = $newParts -join ", "
    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
